$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column I: "Gutenberg" ---
# Header (bold style, matches H1)
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "Gutenberg"

# Row 6 is the only "yes" -- write it first so the shared-string table
# gets "Gutenberg yes" appended before "Gutenberg no".
$ws.Range("H6").Copy()
$ws.Range("I6").PasteSpecial(-4122)
$ws.Range("I6").Value = "Gutenberg yes"

# Remaining data rows (2-5, 7-48) are all "Gutenberg no"
$noRows = 2,3,4,5,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48
foreach ($r in $noRows) {
    $srcCell = "H" + $r
    $dstCell = "I" + $r
    $ws.Range($srcCell).Copy()
    $ws.Range($dstCell).PasteSpecial(-4122)
    $ws.Range($dstCell).Value = "Gutenberg no"
}

# --- D2 typo edit ---
$ws.Range("D2").Value = "Works printedin England, 1784-1807"

# --- View state: zoom + selection ---
$excel.ActiveWindow.Zoom = 59
$ws.Range("D2").Select() | Out-Null
